# Regenerate the save_data sheet's "K" column (column G) values.
# The workbook's column G header is "K" (strike count), and this script
# writes the freshly recalculated K values for each data row (rows 2-73),
# matching the regenerated std/mean + s_vals computation described in the
# commit message. Only column G changes; all other columns/rows are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 0
    4  = 3
    5  = 1
    6  = 0
    7  = 1
    8  = 0
    9  = 1
    10 = 1
    11 = 2
    13 = 3
    14 = 1
    15 = 2
    16 = 0
    17 = 0
    18 = 2
    19 = 3
    20 = 1
    21 = 1
    22 = 3
    23 = 1
    24 = 0
    25 = 1
    26 = 3
    27 = 0
    28 = 0
    30 = 2
    31 = 0
    32 = 0
    33 = 0
    34 = 3
    35 = 1
    36 = 0
    37 = 0
    38 = 3
    39 = 1
    40 = 2
    41 = 1
    42 = 1
    43 = 0
    44 = 1
    45 = 2
    46 = 1
    47 = 0
    48 = 0
    49 = 2
    50 = 0
    51 = 0
    52 = 2
    53 = 0
    54 = 1
    55 = 0
    56 = 0
    57 = 3
    58 = 1
    59 = 0
    61 = 1
    62 = 1
    63 = 0
    64 = 1
    65 = 1
    66 = 0
    68 = 2
    69 = 2
    70 = 1
    71 = 1
    72 = 3
    73 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
